# Scheduled runner update: refresh computed Leve-profit columns (H:N) on
# several rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with
# newer market-board pricing data.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1041.2858
$ws.Range("J19").Value = 1112
$ws.Range("L19").Value = 1112
$ws.Range("N19").Value = -1462
# Row 32
$ws.Range("H32").Value = 800
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 800
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 800
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1452
# Row 64
$ws.Range("H64").Value = 4470
$ws.Range("J64").Value = 4470
$ws.Range("L64").Value = 4470
$ws.Range("N64").Value = -4966
# Row 67
$ws.Range("H67").Value = 4470
$ws.Range("J67").Value = 4470
$ws.Range("L67").Value = 4470
$ws.Range("N67").Value = -6186
# Row 80
$ws.Range("H80").Value = 1333.6364
$ws.Range("I80").Value = 1199.2858
$ws.Range("J80").Value = 1568.75
$ws.Range("K80").Value = 3597.8574
$ws.Range("L80").Value = 4706.25
$ws.Range("M80").Value = -2599.8574
$ws.Range("N80").Value = -6702.25
# Row 83
$ws.Range("H83").Value = 1333.6364
$ws.Range("I83").Value = 1199.2858
$ws.Range("J83").Value = 1568.75
$ws.Range("K83").Value = 10793.5722
$ws.Range("L83").Value = 14118.75
$ws.Range("M83").Value = -5801.572200000001
$ws.Range("N83").Value = -24102.75
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 127
$ws.Range("H127").Value = 2000
$ws.Range("I127").Value = 2000
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 6000
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -1040
$ws.Range("N127").ClearContents()
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 138
$ws.Range("H138").Value = 3134.4348
$ws.Range("I138").Value = 1442.4286
$ws.Range("J138").Value = 3874.6875
$ws.Range("K138").Value = 4327.2858
$ws.Range("L138").Value = 11624.0625
$ws.Range("M138").Value = 812.7142000000003
$ws.Range("N138").Value = -21904.0625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2584.3215
$ws.Range("J32").Value = 8999
$ws.Range("L32").Value = 8999
$ws.Range("N32").Value = -9573
# Row 97
$ws.Range("H97").Value = 625
$ws.Range("I97").Value = 625
$ws.Range("K97").Value = 625
$ws.Range("M97").Value = -129
# Row 132
$ws.Range("H132").Value = 2811
$ws.Range("I132").Value = 2678.889
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8036.667
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5506.667
$ws.Range("N132").Value = -17060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 8068.615
$ws.Range("I134").Value = 8175
$ws.Range("J134").Value = 7898.4
$ws.Range("K134").Value = 24525
$ws.Range("L134").Value = 23695.2
$ws.Range("M134").Value = -21990
$ws.Range("N134").Value = -28765.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1557.3043
$ws.Range("J58").Value = 1497.3846
$ws.Range("L58").Value = 1497.3846
$ws.Range("N58").Value = -1903.3846
# Row 132
$ws.Range("H132").Value = 3168.7144
$ws.Range("I132").Value = 3168.7144
$ws.Range("K132").Value = 9506.143199999999
$ws.Range("M132").Value = -6976.143199999999
# Row 134
$ws.Range("H134").Value = 3945.0833
$ws.Range("J134").Value = 4473.6665
$ws.Range("L134").Value = 13420.9995
$ws.Range("N134").Value = -18490.9995
# Row 136
$ws.Range("H136").Value = 1557.3043
$ws.Range("J136").Value = 1497.3846
$ws.Range("L136").Value = 4492.1538
$ws.Range("N136").Value = -9592.1538

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 1615.6666
$ws.Range("J107").Value = 1599.8
$ws.Range("L107").Value = 4799.4
$ws.Range("N107").Value = -8639.4
# Row 131
$ws.Range("H131").Value = 30
$ws.Range("I131").Value = 30
$ws.Range("K131").Value = 90
$ws.Range("M131").Value = 4950

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3285.8333
$ws.Range("I132").Value = 3285.8333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9857.499899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7327.499899999999
$ws.Range("N132").ClearContents()
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 9749.75
$ws.Range("I17").Value = 8499.5
$ws.Range("J17").Value = 11000
$ws.Range("K17").Value = 8499.5
$ws.Range("L17").Value = 11000
$ws.Range("M17").Value = -8329.5
$ws.Range("N17").Value = -11340
# Row 22
$ws.Range("H22").Value = 1950.0416
$ws.Range("I22").Value = 1956.375
$ws.Range("J22").Value = 1937.375
$ws.Range("K22").Value = 1956.375
$ws.Range("L22").Value = 1937.375
$ws.Range("M22").Value = -1661.375
$ws.Range("N22").Value = -2527.375
# Row 27
$ws.Range("H27").Value = 1950.0416
$ws.Range("I27").Value = 1956.375
$ws.Range("J27").Value = 1937.375
$ws.Range("K27").Value = 1956.375
$ws.Range("L27").Value = 1937.375
$ws.Range("M27").Value = -1849.375
$ws.Range("N27").Value = -2151.375
# Row 63
$ws.Range("H63").Value = 29817
$ws.Range("J63").Value = 29817
$ws.Range("L63").Value = 29817
$ws.Range("N63").Value = -31315
# Row 66
$ws.Range("H66").Value = 29817
$ws.Range("J66").Value = 29817
$ws.Range("L66").Value = 89451
$ws.Range("N66").Value = -96939
# Row 68
$ws.Range("H68").Value = 35699.8
$ws.Range("I68").Value = 3500
$ws.Range("J68").Value = 43749.75
$ws.Range("K68").Value = 3500
$ws.Range("L68").Value = 43749.75
$ws.Range("M68").Value = -2751
$ws.Range("N68").Value = -45247.75
# Row 71
$ws.Range("H71").Value = 35699.8
$ws.Range("I71").Value = 3500
$ws.Range("J71").Value = 43749.75
$ws.Range("K71").Value = 17500
$ws.Range("L71").Value = 218748.75
$ws.Range("M71").Value = -13756
$ws.Range("N71").Value = -226236.75
# Row 74
$ws.Range("H74").Value = 36248.5
$ws.Range("I74").Value = 22500
$ws.Range("J74").Value = 49997
$ws.Range("K74").Value = 22500
$ws.Range("L74").Value = 49997
$ws.Range("M74").Value = -21502
$ws.Range("N74").Value = -51993
# Row 77
$ws.Range("H77").Value = 36248.5
$ws.Range("I77").Value = 22500
$ws.Range("J77").Value = 49997
$ws.Range("K77").Value = 67500
$ws.Range("L77").Value = 149991
$ws.Range("M77").Value = -62508
$ws.Range("N77").Value = -159975

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3966.5
$ws.Range("I62").Value = 3966.5
$ws.Range("K62").Value = 3966.5
$ws.Range("M62").Value = -3342.5
# Row 65
$ws.Range("H65").Value = 3966.5
$ws.Range("I65").Value = 3966.5
$ws.Range("K65").Value = 19832.5
$ws.Range("M65").Value = -16712.5
# Row 132
$ws.Range("H132").Value = 1326.1666
$ws.Range("I132").Value = 1310.8
$ws.Range("K132").Value = 3932.4
$ws.Range("M132").Value = -1402.4
# Row 136
$ws.Range("H136").Value = 3548.1853
$ws.Range("I136").Value = 3412.389
$ws.Range("J136").Value = 3819.7778
$ws.Range("K136").Value = 10237.167
$ws.Range("L136").Value = 11459.3334
$ws.Range("M136").Value = -7687.167000000001
$ws.Range("N136").Value = -16559.3334
